# TimeManage_ZYJ.xlsx - add the 2012.4.10 "git" progress entry as a new row
# (row 6) in Sheet1, mirroring the formatting of the row above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 / column A: date-like text "2012.4.10" -----------------------
# A plain Value assignment of a date-shaped string ("2012.4.10") gets
# auto-parsed into a date serial by Excel's smart typing, which is not
# what the source file has (it is stored as literal text). Build it as a
# text formula result in a scratch cell, then paste-special just the
# *value* into A6 so it lands as literal text instead of a date number.
$ws.Range("F1").Formula = "=""2012.4.10"""
$ws.Range("F1").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$ws.Range("F1").ClearContents() | Out-Null

# --- Copy the formatting of row 4 (A4:B4) down onto row 6 ---------------
# This mirrors the style used by the other single-line rows (A6 like A4,
# B6 like B4) without disturbing the values we just set.
$ws.Range("A4:B4").Copy() | Out-Null
$ws.Range("A6:B6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Fill in the rest of row 6 -------------------------------------------
$ws.Range("B6").Value = "学习git相关操作，对git的分支操作有了更多的了解。上传git搭建和使用文档"
$ws.Range("D6").Value = 2

# Row 6 is a bit taller than the default to fit its content.
$ws.Rows("6").RowHeight = 27

# Leave the selection where the author ended up after entering the row.
$ws.Range("F5").Select() | Out-Null
